$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167127132415771
$ws.Range("B1").Value = 5.207211017608643
$ws.Range("C1").Value = 1.529977321624756
$ws.Range("D1").Value = 1.008427143096924
$ws.Range("E1").Value = 0.5507215857505798
